$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4 content: AutoSPInstaller Build 3.96 (SP2013 SP1 Managed Accounts guidance)
# Shared-string insertion order matters (matches original commit's sharedStrings.xml
# ordering): Configuration Section (B4) -> Examples (D4) -> Notes (C4).

$bText = " <ManagedAccounts>"

$dText = @'
<ManagedAccount CommonName="spservice">
                <Username>DOMAIN\SP_Services</Username>
                <Password></Password>
            </ManagedAccount>
            <ManagedAccount CommonName="Portal">
                <Username>DOMAIN\SP_PortalAppPool</Username>
                <Password></Password>
            </ManagedAccount>
            <ManagedAccount CommonName="MySiteHost">
                <Username>DOMAIN\SP_ProfilesAppPool</Username>
                <Password></Password>
            </ManagedAccount>
            <ManagedAccount CommonName="SearchService">
                <Username>DOMAIN\SP_SearchService</Username>
                <Password></Password>
            </ManagedAccount>
'@

$cText = "When provisioning Managed Accounts it is important to leave the 'CommonName' property set to what they are for the 4 default accounts. This property is now bound to the default 'Portal' and 'My Site' web applications; and the Search Service Applications - so should be left in place."

# Column A: Configuration Version
$ws.Range("A4").Value = 3.96
$ws.Range("A4").VerticalAlignment = -4160

# Column B: Configuration Section (same alignment style as A2/A3 - vertical top, no wrap)
$ws.Range("B4").Value = $bText
$ws.Range("B4").VerticalAlignment = -4160

# Column D: Examples (vertical top + wrap, same style as B2/B3/C2/C3/D2/D3)
$ws.Range("D4").Value = $dText
$ws.Range("D4").VerticalAlignment = -4160
$ws.Range("D4").WrapText = $true

# Column C: Notes (vertical top + wrap)
$ws.Range("C4").Value = $cText
$ws.Range("C4").VerticalAlignment = -4160
$ws.Range("C4").WrapText = $true

$ws.Rows.Item(4).RowHeight = 409.5

$ws.Range("A4").Select()
